$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Feuil1" - the active sheet

# ---------------------------------------------------------------------------
# 1) Row 27 ("Press start"): mark as fully done (green, 100%)
#    B27 gets the "done" look (same look already used e.g. by B62),
#    E27 gets the "100%, blue font" percentage look (same as E9) and value 1.
# ---------------------------------------------------------------------------
$ws.Range("B62").Copy()
$ws.Range("B27").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E9").Copy()
$ws.Range("E27").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E27").Value2 = 1

# ---------------------------------------------------------------------------
# 2) Row 28 ("Mode screen"): same treatment as row 27
# ---------------------------------------------------------------------------
$ws.Range("B62").Copy()
$ws.Range("B28").PasteSpecial(-4122)
$ws.Range("E9").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value2 = 1

# ---------------------------------------------------------------------------
# 3) Row 46 ("Songwheel screen"): song selection pack is in progress (20%)
#    B46 takes the light-blue header look already used on F1,
#    E46 takes the red percentage look already used on E64/E63(before),
#    F46 gets a new comment "Pack ok, song en cours".
# ---------------------------------------------------------------------------
$ws.Range("F1").Copy()
$ws.Range("B46").PasteSpecial(-4122)
$ws.Range("E64").Copy()
$ws.Range("E46").PasteSpecial(-4122)
$ws.Range("E46").Value2 = 0.2
$ws.Range("F46").Value2 = "Pack ok, song en cours"

# ---------------------------------------------------------------------------
# 4) Row 63 ("FadeManager"): now 90% done, with a comment on what remains
#    E63 takes the red percentage look (same as E64), value 0.9,
#    F63 gets a new comment "D'autres trucs à rajouter par la suite ?".
#    B63 keeps its own distinctive style, only its fill becomes the lighter
#    blue-grey tone already used elsewhere (e.g. on B46/F1).
# ---------------------------------------------------------------------------
$ws.Range("E64").Copy()
$ws.Range("E63").PasteSpecial(-4122)
$ws.Range("E63").Value2 = 0.9
$ws.Range("F63").Value2 = "D'autres trucs à rajouter par la suite ?"
$ws.Range("B63").Interior.Color = $ws.Range("F1").Interior.Color

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 5) Update the saved scroll position / selection of the sheet
# ---------------------------------------------------------------------------
$ws.Range("D63").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 51
$win.ScrollColumn = 1
